$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 176
$ws.Range("J2").Value = 633
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 164
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 109
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 69
$ws.Range("T2").Value = 102
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 972
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 900
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 9
